# Generate Report for Handback
#
# The handback has completed: both locales (zh-cn, de-de) are now in sync
# with en-US. Refresh the Overview status column, stamp the handback
# datetime for each locale, and record the "Latest Target File" /
# "Latest Handback File" columns (with their hyperlinks) for every row of
# the per-locale report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: both files are now handed back & in sync with en-US.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn report sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Latest Handback DateTime for this locale.
$zhcn.Range("H2").Value = "2016-03-23 20:54:13"
$zhcn.Range("H3").Value = "2016-03-23 20:54:13"

# Row 2 - 4f3491b5-b04d-497b-a309-93df6e5fb3f9
$zhcn.Range("F2").Value = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/3c8ec4a0769d5d75b18d5368f26c3814904a30c4/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md", "", "", "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md") | Out-Null

$zhcn.Range("G2").Value = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3869c591d223b0be9e884b4e91550e1de6af541e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf", "", "", "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf") | Out-Null

# Row 3 - a53b297b-d811-421c-9e68-f34339466385
$zhcn.Range("F3").Value = "a53b297b-d811-421c-9e68-f34339466385.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/3c8ec4a0769d5d75b18d5368f26c3814904a30c4/e2e/a53b297b-d811-421c-9e68-f34339466385.md", "", "", "a53b297b-d811-421c-9e68-f34339466385.md") | Out-Null

$zhcn.Range("G3").Value = "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3869c591d223b0be9e884b4e91550e1de6af541e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf", "", "", "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de report sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Latest Handback DateTime for this locale.
$dede.Range("H2").Value = "2016-03-23 20:54:22"
$dede.Range("H3").Value = "2016-03-23 20:54:22"

# Row 2 - 4f3491b5-b04d-497b-a309-93df6e5fb3f9
$dede.Range("F2").Value = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/3c8ec4a0769d5d75b18d5368f26c3814904a30c4/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md", "", "", "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md") | Out-Null

$dede.Range("G2").Value = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e433d3a6fc9800cfc40825cd21353444ed62014/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf", "", "", "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf") | Out-Null

# Row 3 - a53b297b-d811-421c-9e68-f34339466385
$dede.Range("F3").Value = "a53b297b-d811-421c-9e68-f34339466385.md"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/3c8ec4a0769d5d75b18d5368f26c3814904a30c4/e2e/a53b297b-d811-421c-9e68-f34339466385.md", "", "", "a53b297b-d811-421c-9e68-f34339466385.md") | Out-Null

$dede.Range("G3").Value = "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e433d3a6fc9800cfc40825cd21353444ed62014/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf", "", "", "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf") | Out-Null
